# This script reorders the per-row "measurement" data (Fecha, Volumen,
# Precio minimo, Precio maximo, Precio promedio ponderado, Precio $/Kg)
# across rows 2-10 and 12-20 of the active sheet, according to a fixed
# permutation. Row 11 is left untouched. All other columns (A,B,C,E,F,G,
# H,I,N,O,Q,R) are identical across these rows so they do not need to be
# touched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# columns that carry the data which gets shuffled between rows
$cols = @("D", "J", "K", "L", "M", "P")

# mapping: target row -> source row (values currently sitting in source
# row must end up in target row)
$mapping = @{
    2  = 16
    3  = 10
    4  = 14
    5  = 2
    6  = 18
    7  = 9
    8  = 4
    9  = 15
    10 = 5
    12 = 6
    13 = 3
    14 = 8
    15 = 19
    16 = 7
    17 = 12
    18 = 20
    19 = 13
    20 = 17
}

# 1) Snapshot the current values for every source row/column before any
#    writes happen, so overwriting one row doesn't clobber data that is
#    still needed for another target row.
$snapshot = @{}
foreach ($row in $mapping.Values) {
    if (-not $snapshot.ContainsKey($row)) {
        $rowVals = @{}
        foreach ($col in $cols) {
            $rowVals[$col] = $ws.Range("$col$row").Value()
        }
        $snapshot[$row] = $rowVals
    }
}

# 2) Apply the snapshot values to the target rows.
foreach ($targetRow in $mapping.Keys) {
    $sourceRow = $mapping[$targetRow]
    $rowVals = $snapshot[$sourceRow]
    foreach ($col in $cols) {
        $ws.Range("$col$targetRow").Value = $rowVals[$col]
    }
}
